# "Loan RBI, Variable Instalments"
# Switch focus to the "Repayment schedule" sheet and insert a new
# (currently blank) column before column N, shifting the old N/O/P
# columns ("Late" / heading / "Outstanding") one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make "Repayment schedule" the active sheet/tab (was "Edit Repayment
# Schedule"); this also clears the stale tabSelected flag on that sheet.
$ws.Activate()

# Insert a blank column at N, pushing the existing N:P columns to O:Q.
# A freshly inserted column inherits the width of the column to its
# left (M), so mirror that explicitly.
$newColWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $newColWidth

# Leave the cursor parked on the new rightmost data column.
$ws.Range("S6").Select() | Out-Null
